# Commit: "Update boh di file"
#
# 1) The cached "datetimeFigureOut" field text on the slide master, every
#    slide layout, and the notes master is bumped from 10/04/25 -> 11/04/25
#    (the deck was simply reopened/resaved a day later, which refreshes the
#    auto date placeholders).
# 2) Slide 2's title run "Dataset " + "refresh" is collapsed into a single
#    run reading just "Dataset".

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Segnaposto data*" -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "10/04/25") {
                $tr.Text = "11/04/25"
            }
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout hanging off the master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Notes master
Update-DatePlaceholder $p.NotesMaster.Shapes

# Slide 2 title: "Dataset refresh" -> "Dataset"
$slide2 = $p.Slides.Item(2)
for ($i = 1; $i -le $slide2.Shapes.Count; $i++) {
    $shp = $slide2.Shapes.Item($i)
    if ($shp.Name -eq "Titolo 1" -and $shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq "Dataset refresh") {
            $shp.TextFrame.TextRange.Text = "Dataset"
        }
    }
}
